# Apply weekly price-data refresh to the Jengibre (ginger) sheet.
# The data rows (2-19) keep their market/variety/quality metadata but several
# rows get updated Fecha (D), Volumen (J), Precio minimo/maximo/promedio
# (K/L/M) and Precio $/Kg (P) values, matching the new weekly source extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = 44895
    "J2"  = 30
    "K2"  = 18000
    "L2"  = 18000
    "M2"  = 18000
    "P2"  = 1385

    "D3"  = 44797
    "J3"  = 60

    "D4"  = 44915
    "J4"  = 50
    "K4"  = 18000
    "L4"  = 18000
    "M4"  = 18000
    "P4"  = 1385

    "D5"  = 44868
    "J5"  = 30
    "K5"  = 18000
    "L5"  = 18000
    "M5"  = 18000
    "P5"  = 1385

    "D6"  = 44880
    "K6"  = 17000
    "L6"  = 17000
    "M6"  = 17000
    "P6"  = 1308

    "D7"  = 44839
    "J7"  = 40
    "K7"  = 15000
    "L7"  = 16000
    "M7"  = 15500
    "P7"  = 1192

    "D8"  = 44804
    "J8"  = 40
    "K8"  = 12000
    "L8"  = 13000
    "M8"  = 12500
    "P8"  = 962

    "D9"  = 44922
    "J9"  = 30
    "L9"  = 17000
    "M9"  = 17000
    "P9"  = 1308

    "D10" = 44810
    "J10" = 40

    "D12" = 44959
    "K12" = 19000
    "L12" = 19000
    "M12" = 19000
    "P12" = 1462

    "D13" = 44943
    "K13" = 17000
    "L13" = 17000
    "M13" = 17000
    "P13" = 1308

    "D14" = 44930

    "D15" = 44859
    "K15" = 13000
    "L15" = 13000
    "M15" = 13000
    "P15" = 1000

    "D16" = 44841

    "D17" = 44894

    "D18" = 44832
    "J18" = 60
    "L18" = 18000
    "M18" = 17500
    "P18" = 1346

    "D19" = 44846
    "K19" = 18000
    "L19" = 18000
    "M19" = 18000
    "P19" = 1385
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
